# Add the "2022-Q1" fund-holdings sheet (inserted right before the "总计"
# summary sheet) and extend the "总计" sheet with a new leading row for it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet, positioned immediately before "总计"
# ---------------------------------------------------------------------------
# NOTE: the sheet object passed as the "Before" argument to Worksheets.Add()
# gets its underlying reference repointed to the freshly-created sheet once
# Add() returns, so it must never be reused afterwards. Pass a throwaway
# Item(...) lookup inline instead of caching "总计" in a variable.
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Copy the A1:H25 block (header + same styling) from the most recent quarter
# sheet so the new sheet starts with identical formatting/column layout.
$wb.Worksheets.Item("2021-Q4").Range("A1:H25").Copy($newSheet.Range("A1"))

# Fund-holdings rows for 2022-Q1: index, code, name, fund scale, stock
# position, position ratio, held value (亿元), position rank.
$data2022Q1 = @(
  ,@(0, '161040', '富国创业板两年定期开放混合', '35.36', '83.63', '3.22', '1.1386', 10)
  ,@(1, '160916', '大成优选混合(LOF)', '16.14', '89.35', '3.98', '0.6424', 8)
  ,@(2, '010738', '大成优选升级一年持有期混合A', '3.79', '89.02', '7.18', '0.2721', 6)
  ,@(3, '009914', '富国成长动力混合', '8.61', '82.19', '2.47', '0.2127', 10)
  ,@(4, '003131', '国寿安保强国智造灵活配置混合', '5.76', '86.19', '3.64', '0.2097', 2)
  ,@(5, '005683', '国寿安保华兴灵活配置混合', '3.43', '89.19', '3.38', '0.1159', 2)
  ,@(6, '004760', '国寿安保稳瑞混合A', '8.20', '21.27', '0.90', '0.0738', 4)
  ,@(7, '000969', '前海开源大安全核心精选灵活配置混合', '1.39', '91.04', '4.21', '0.0585', 5)
  ,@(8, '003165', '鹏华弘嘉灵活配置混合A', '1.53', '93.95', '2.74', '0.0419', 10)
  ,@(9, '001060', '前海开源高端装备制造灵活配置混合', '0.97', '89.88', '4.28', '0.0415', 5)
  ,@(10, '011734', '国寿安保裕丰混合型证券投资基金A', '5.01', '20.33', '0.80', '0.0401', 1)
  ,@(11, '010205', '国寿安保裕安混合A', '3.44', '29.13', '1.08', '0.0372', 3)
  ,@(12, '011484', '申万菱信宜选混合A', '4.30', '22.45', '0.67', '0.0288', 10)
  ,@(13, '004761', '国寿安保稳瑞混合C', '2.52', '21.27', '0.90', '0.0227', 4)
  ,@(14, '009128', '明亚价值长青混合A', '0.38', '49.48', '4.06', '0.0154', 4)
  ,@(15, '003166', '鹏华弘嘉灵活配置混合C', '0.56', '93.95', '2.74', '0.0153', 10)
  ,@(16, '010765', '国寿安保华丰混合A', '0.42', '88.40', '3.08', '0.0129', 3)
  ,@(17, '011735', '国寿安保裕丰混合型证券投资基金C', '1.37', '20.33', '0.80', '0.0110', 1)
  ,@(18, '011485', '申万菱信宜选混合C', '1.07', '22.45', '0.67', '0.0072', 10)
  ,@(19, '010739', '大成优选升级一年持有期混合C', '0.09', '89.02', '7.18', '0.0065', 6)
  ,@(20, '010206', '国寿安保裕安混合C', '0.20', '29.13', '1.08', '0.0022', 3)
  ,@(21, '005146', '兴银丰润灵活配置混合', '0.05', '93.36', '3.29', '0.0016', 6)
  ,@(22, '010766', '国寿安保华丰混合C', '0.00', '88.40', '3.08', 0, 3)
  ,@(23, '009129', '明亚价值长青混合C', '0.00', '49.48', '4.06', 0, 4)
)

for ($i = 0; $i -lt $data2022Q1.Count; $i++) {
  $r = $i + 2
  $row = $data2022Q1[$i]
  $newSheet.Cells.Item($r, 1).Value = $row[0]
  $newSheet.Cells.Item($r, 2).Value = $row[1]
  $newSheet.Cells.Item($r, 3).Value = $row[2]
  $newSheet.Cells.Item($r, 4).Value = $row[3]
  $newSheet.Cells.Item($r, 5).Value = $row[4]
  $newSheet.Cells.Item($r, 6).Value = $row[5]
  $newSheet.Cells.Item($r, 7).Value = $row[6]
  $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2) Extend the "总计" sheet with a new first data row for 2022-Q1, shifting
#    the existing quarters down by one row.
# ---------------------------------------------------------------------------
# Re-fetch "总计" fresh (now sitting after "2022-Q1" in the tab order) rather
# than reusing any earlier reference.
$total = $wb.Worksheets.Item("总计")

$totalsData = @(
  ,@(0, '2022-Q1', 24, 3.01)
  ,@(1, '2021-Q4', 40, 6.63)
  ,@(2, '2021-Q3', 27, 3.49)
  ,@(3, '2021-Q2', 16, 1.58)
  ,@(4, '2021-Q1', 17, 5.04)
  ,@(5, '2020-Q4', 15, 4.51)
)

# Seed formatting for the brand-new last row (row 7) by copying column A's
# styled cell down from the previous last row before the values are written.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

for ($i = 0; $i -lt $totalsData.Count; $i++) {
  $r = $i + 2
  $row = $totalsData[$i]
  $total.Cells.Item($r, 1).Value = $row[0]
  $total.Cells.Item($r, 2).Value = $row[1]
  $total.Cells.Item($r, 3).Value = $row[2]
  $total.Cells.Item($r, 4).Value = $row[3]
}
